$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.711.17"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.633.41"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'217.94"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").Value = "'18.93"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("D11").Value = "'0.0841"
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").Value = "1.626.18"
$ws.Range("E13").Value = "  -1.64%  "
$ws.Range("E14").Value = "  -2.82%  "
$ws.Range("E15").Value = "  -2.30%  "
$ws.Range("D16").Value = "'64.01"
$ws.Range("E16").Value = "  -2.73%  "
$ws.Range("D17").Value = "26.683.89"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "0.0₃0724"
$ws.Range("E18").Value = "  -3.11%  "
$ws.Range("D19").Value = "'211.09"
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("E23").Value = "  -4.32%  "
$ws.Range("E24").Value = "  -3.26%  "
$ws.Range("D25").Value = "'146.72"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("D28").Value = "'6.98"
$ws.Range("E28").Value = "  -3.29%  "
$ws.Range("D29").Value = "'15.53"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("D30").Value = "'0.0499"
$ws.Range("E30").Value = "  -4.41%  "
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "'2.95"
$ws.Range("E33").Value = "  -2.77%  "
$ws.Range("D34").Value = "1.261.38"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("D38").Value = "'0.525"
$ws.Range("E38").Value = "  -3.41%  "
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("E40").Value = "  -4.06%  "
$ws.Range("D41").Value = "'0.797"
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("E42").Value = "  -3.21%  "
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "1.772.00"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.15"
$ws.Range("E44").Value = "  -4.54%  "
$ws.Range("D45").Value = "'91.42"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").Value = "'59.70"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "'1.56"
$ws.Range("E47").Value = "  -4.31%  "
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("E51").Value = "  -2.73%  "
